$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-533) holds a "Förändrad" (changed) date that needs to be
# bumped from 2023-09-19 (serial 45188) to 2023-09-20 (serial 45189) for every
# data row on the sheet.
$ws.Range("C2:C533").Value = 45189
